$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D3").Value = "2016-01-08 13:35:13"
$wsZh.Range("G3").Value = "2016-01-08 13:36:00"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D3").Value = "2016-01-08 13:35:26"
$wsDe.Range("G3").Value = "2016-01-08 13:36:22"
